$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price / Volume(1h) columns) per latest scrape.
# D-column ("Price") values look numeric (e.g. "71.223.13", "0.0000258") but
# must remain literal text, exactly as scraped, so we briefly force a text
# number format while assigning them, then restore the default cell style
# so no formatting residue is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.223.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.848.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "696.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.845.51"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.495.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.847.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.234.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.47%  "
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.002.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.64%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  +2.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.799.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +3.19%  "
$ws.Range("E40").Value = "  +12.02%  "
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("E43").Value = "  +5.80%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000305"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "419.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.17%  "
$ws.Range("E51").Value = "  +1.01%  "
